$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costs")

# H3 and H4 no longer derive their "Marginal cost" from VOM (column D);
# they are now hard-coded to 0. This breaks them out of the H3:H5 shared
# formula group, leaving H5 (=D5) as the sole remaining formula there.
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0

# Move / restore the active selection to H5 (last edited cell).
$ws.Range("H5").Select()
